# A new weekly price record for "Valencia" oranges (Provincia de Melipilla,
# $/caja 15 kilos granel) needs to be inserted at the top of the data block
# (row 193), pushing the existing rows 193-239 down to 194-240.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 193; Excel shifts rows 193:239 down to 194:240 and
# carries the row-193 formatting (e.g. the date style on column D) along.
$ws.Rows("193:193").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(193, 1).Value  = 11
$ws.Cells.Item(193, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(193, 3).Value  = "Bíobío"
$ws.Cells.Item(193, 4).Value  = 44642
$ws.Cells.Item(193, 5).Value  = 8
$ws.Cells.Item(193, 6).Value  = "Fruta"
$ws.Cells.Item(193, 7).Value  = 100102
$ws.Cells.Item(193, 8).Value  = "Cítricos"
$ws.Cells.Item(193, 9).Value  = 100102005
$ws.Cells.Item(193, 10).Value = "Naranja"
$ws.Cells.Item(193, 11).Value = "Valencia"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 220
$ws.Cells.Item(193, 14).Value = 9000
$ws.Cells.Item(193, 15).Value = 10000
$ws.Cells.Item(193, 16).Value = 9545
$ws.Cells.Item(193, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(193, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(193, 19).Value = 636
$ws.Cells.Item(193, 20).Value = 15
